# Auto commit at 2025-08-31  8:09:11.14
# Updates metric figures (charging/income totals) that cascade into the
# "today" summary sheet, appends a 2026 row to the year-trend sheets, and
# nudges the remembered cell selection on the touched sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metrics sheet: refresh the monthly/yearly/cumulative figures.
# ---------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 507971.08
$metrics.Range("B3").Value = 435488.04000000004
$metrics.Range("B4").Value = 160782.65
$metrics.Range("B5").Value = 19863
$metrics.Range("B6").Value = 3903599.6499999994
$metrics.Range("B7").Value = 3314202.6999999997
$metrics.Range("B8").Value = 1124425.21
$metrics.Range("B9").Value = 150551
$metrics.Range("B10").Value = 32368923.45099983
$metrics.Range("B11").Value = 19344072.770000003
$metrics.Range("B12").Value = 11406134.100000001
$metrics.Range("B13").Value = 1248178

$metrics.Activate() | Out-Null
$metrics.Range("G38").Select() | Out-Null

# ---------------------------------------------------------------------
# today sheet: keep the TODAY()-1 formula in A1 as-is, just move the
# remembered selection (the cached date recalculates automatically).
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("G9").Select() | Out-Null

# ---------------------------------------------------------------------
# csdjzqs sheet: give column C the same number format as column B, then
# append the 2026 row.
# ---------------------------------------------------------------------
$csdjzqs = $wb.Worksheets.Item("csdjzqs")
$ndzsrqs = $wb.Worksheets.Item("ndzsrqs")

$csdjzqs.Range("C2:C9").NumberFormat = $ndzsrqs.Range("B2").NumberFormat

$csdjzqs.Range("A10").Value = "2026年"
$csdjzqs.Range("B10").Value = 123456
$csdjzqs.Range("C10").Value = 12345
$csdjzqs.Range("B10").NumberFormat = $csdjzqs.Range("B9").NumberFormat
$csdjzqs.Range("C10").NumberFormat = $ndzsrqs.Range("B2").NumberFormat

$csdjzqs.Activate() | Out-Null
$csdjzqs.Range("E20").Select() | Out-Null

# ---------------------------------------------------------------------
# ndzsrqs sheet: append the 2026 row.
# ---------------------------------------------------------------------
$ndzsrqs.Range("A10").Value = 2026
$ndzsrqs.Range("B10").Value = 123456
$ndzsrqs.Range("B10").NumberFormat = $ndzsrqs.Range("B9").NumberFormat

$ndzsrqs.Activate() | Out-Null
$ndzsrqs.Range("C12").Select() | Out-Null

# ---------------------------------------------------------------------
# Restore the workbook's originally active tab (bksr) so the document-
# level active sheet/bookview stays unchanged, matching the diff (only
# the per-sheet remembered selections move).
# ---------------------------------------------------------------------
$bksr = $wb.Worksheets.Item("bksr")
$bksr.Activate() | Out-Null
